$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I: "H07" (Homework 07)
$ws.Range("I1").Value = "H07"

# Homework 07 scores for each student (rows 2-16)
$ws.Range("I2").Value = 9.5
$ws.Range("I3").Value = 10
$ws.Range("I4").Value = 9
$ws.Range("I5").Value = 9
$ws.Range("I6").Value = 10
$ws.Range("I7").Value = 8.5
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 9.5
$ws.Range("I10").Value = 8.5

# Rows 11-16 also need center alignment to match the rest of the column
$ws.Range("I11").Value = 9.5
$ws.Range("I11").HorizontalAlignment = -4108

$ws.Range("I12").Value = 9.5
$ws.Range("I12").HorizontalAlignment = -4108

# Row 13 is given half credit via a formula, highlighted like D15
$ws.Range("I13").Formula = "=9.5/2"
$ws.Range("I13").HorizontalAlignment = -4108
$ws.Range("I13").Interior.Color = 65535

$ws.Range("I14").Value = 9
$ws.Range("I14").HorizontalAlignment = -4108

$ws.Range("I15").Value = 9.5
$ws.Range("I15").HorizontalAlignment = -4108

$ws.Range("I16").Value = 9.5
$ws.Range("I16").HorizontalAlignment = -4108

# Totals row (17) - keep blank, but align like the rest of the column
$ws.Range("I17").HorizontalAlignment = -4108

# Move the active selection to B34
$ws.Range("B34").Select() | Out-Null
